$p = $ppt.ActivePresentation

# Insert a new slide ("Title and Content" layout = CustomLayout index 2)
# right before the current last slide ("Исходники и материалы"), which sits
# at index 10. The new slide becomes slide #10 and the old last slide is
# pushed down to slide #11.
$newSlide = $p.Slides.Add(10, 2)

# Title placeholder -> "Источники"
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Источники"

# Body / content placeholder -> bibliography list
$body = $newSlide.Shapes.Item(2)

$body.Left = 60.991338582677166
$body.Top = 213.87181102362206
$body.Width = 756.6781889763779
$body.Height = 179.3107874015748

$tr = $body.TextFrame.TextRange
$tr.Text = "Amara "
$tr.InsertAfter("Graps") | Out-Null
$tr.InsertAfter(", ") | Out-Null
$tr.InsertAfter([char]0x00AB) | Out-Null
$tr.InsertAfter("An Introduction to Wavelets") | Out-Null
$tr.InsertAfter([char]0x00BB) | Out-Null

$tr.InsertAfter("`rRobi") | Out-Null
$tr.InsertAfter(" ") | Out-Null
$tr.InsertAfter("Polikar") | Out-Null
$tr.InsertAfter(", ") | Out-Null
$tr.InsertAfter([char]0x00AB) | Out-Null
$tr.InsertAfter("The Wavelet Tutorial") | Out-Null
$tr.InsertAfter([char]0x00BB) | Out-Null

$tr.InsertAfter("`rMiKXMan") | Out-Null
$tr.InsertAfter(", ") | Out-Null
$tr.InsertAfter([char]0x00AB + "Непрерывное ") | Out-Null
$tr.InsertAfter("wavelet ") | Out-Null
$tr.InsertAfter("преобразование" + [char]0x00BB) | Out-Null

$tr.InsertAfter("`r") | Out-Null
